$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (A: 15.42578125 -> 16.42578125, B: 14.7109375 -> 15.7109375)
# Note: the ColumnWidth COM property is offset from the raw OOXML column
# width (stored width = ColumnWidth + 0.8333333333333334) and gets
# quantized to the nearest 1/6 when Excel re-serializes it, so we pick the
# ColumnWidth value whose serialized width is closest to the target.
$ws.Columns.Item(1).ColumnWidth = 15.592447916666666
$ws.Columns.Item(2).ColumnWidth = 14.877604166666666

# Update cell values in columns A and B for rows 1-32
$ws.Cells.Item(1, 1).Value = -0.054833169165874551
$ws.Cells.Item(1, 2).Value = 0.054583645564107997
$ws.Cells.Item(2, 1).Value = -0.0084738997451658804
$ws.Cells.Item(2, 2).Value = 0.007447390774299123
$ws.Cells.Item(3, 1).Value = 0.095486375985778693
$ws.Cells.Item(3, 2).Value = -0.096096683973755148
$ws.Cells.Item(4, 1).Value = -0.1878959114650911
$ws.Cells.Item(4, 2).Value = 0.18668514358022392
$ws.Cells.Item(5, 1).Value = -0.18068514398232516
$ws.Cells.Item(5, 2).Value = 0.17823223215193185
$ws.Cells.Item(6, 1).Value = -0.056659559603742249
$ws.Cells.Item(6, 2).Value = 0.056616502633896726
$ws.Cells.Item(7, 1).Value = -0.036616503128758637
$ws.Cells.Item(7, 2).Value = 0.036552234274708439
$ws.Cells.Item(8, 1).Value = 0.0056171993817821786
$ws.Cells.Item(8, 2).Value = -0.0056283713340583574
$ws.Cells.Item(9, 1).Value = 0.011628370913252084
$ws.Cells.Item(9, 2).Value = -0.011650283829903962
$ws.Cells.Item(10, 1).Value = 0.017650283410056034
$ws.Cells.Item(10, 2).Value = -0.017653108370744519
$ws.Cells.Item(11, 1).Value = -0.0062828564772061668
$ws.Cells.Item(11, 2).Value = 0.006282350432893935
$ws.Cells.Item(12, 1).Value = -0.00028235085287819928
$ws.Cells.Item(12, 2).Value = 0.00028043314471748459
$ws.Cells.Item(13, 1).Value = 0.0057195664351601394
$ws.Cells.Item(13, 2).Value = -0.0057205409874354629
$ws.Cells.Item(14, 1).Value = -0.027082065063541272
$ws.Cells.Item(14, 2).Value = 0.02705105068045377
$ws.Cells.Item(15, 1).Value = -0.021051051101878215
$ws.Cells.Item(15, 2).Value = 0.021026645461807547
$ws.Cells.Item(16, 1).Value = -0.015026645884658407
$ws.Cells.Item(16, 2).Value = 0.015004234449906839
$ws.Cells.Item(17, 1).Value = -0.0090042348746166567
$ws.Cells.Item(17, 2).Value = 0.0089999995585055004
$ws.Cells.Item(18, 1).Value = -0.036109746238103924
$ws.Cells.Item(18, 2).Value = 0.036096537368983661
$ws.Cells.Item(19, 1).Value = -0.027096537783018793
$ws.Cells.Item(19, 2).Value = 0.02701362552503328
$ws.Cells.Item(20, 1).Value = -0.018013625942719713
$ws.Cells.Item(20, 2).Value = 0.018004259317095261
$ws.Cells.Item(21, 1).Value = -0.0090042597353647835
$ws.Cells.Item(21, 2).Value = 0.0089999995813174749
$ws.Cells.Item(22, 1).Value = -0.093933767174577909
$ws.Cells.Item(22, 2).Value = 0.093625303288030537
$ws.Cells.Item(23, 1).Value = -0.084625303706252986
$ws.Cells.Item(23, 2).Value = 0.084124975021502735
$ws.Cells.Item(24, 1).Value = -0.042124975625957184
$ws.Cells.Item(24, 2).Value = 0.041999999392253073
$ws.Cells.Item(25, 1).Value = -0.15545039666389826
$ws.Cells.Item(25, 2).Value = 0.15481403330498722
$ws.Cells.Item(26, 1).Value = -0.088658993310094303
$ws.Cells.Item(26, 2).Value = 0.088346738410496783
$ws.Cells.Item(27, 1).Value = -0.082346738833397382
$ws.Cells.Item(27, 2).Value = 0.081277997433697902
$ws.Cells.Item(28, 1).Value = -0.075277997865741852
$ws.Cells.Item(28, 2).Value = 0.074537484974597668
$ws.Cells.Item(29, 1).Value = -0.062537485445767871
$ws.Cells.Item(29, 2).Value = 0.062170368188096958
$ws.Cells.Item(30, 1).Value = -0.042170368705961803
$ws.Cells.Item(30, 2).Value = 0.042019879494294088
$ws.Cells.Item(31, 1).Value = -0.027019879989502726
$ws.Cells.Item(31, 2).Value = 0.027000781237669358
$ws.Cells.Item(32, 1).Value = -0.0060007817662874885
$ws.Cells.Item(32, 2).Value = 0.0059999995528823291
